$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 16-21) is being re-sorted / updated so that records are
# grouped by "Periodo Mora" (1803 first, then 1804) instead of by worker.
# New content of rows 16-21, columns C (doc#), D (name), E (periodo),
# F (valor mora), G (salario basico):

$ws.Range("C16").Value = "30873862"
$ws.Range("D16").Value = "KARINA DEL CARMEN VILLALBA BORJA"
$ws.Range("E16").Value = "1803"
$ws.Range("F16").Value = 19200
$ws.Range("G16").Value = 900000

$ws.Range("C17").Value = "1143379280"
$ws.Range("D17").Value = "PAOLA PATRICIA AGUILAR VELASCO"
$ws.Range("E17").Value = "1803"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 781242

$ws.Range("C18").Value = "30873862"
$ws.Range("D18").Value = "KARINA DEL CARMEN VILLALBA BORJA"
$ws.Range("E18").Value = "1804"
$ws.Range("F18").Value = 36000
$ws.Range("G18").Value = 900000

$ws.Range("C19").Value = "1143379280"
$ws.Range("D19").Value = "PAOLA PATRICIA AGUILAR VELASCO"
$ws.Range("E19").Value = "1804"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

$ws.Range("C20").Value = "1143404861"
$ws.Range("D20").Value = "DARIANY CANO DIAZ"
$ws.Range("E20").Value = "1804"
$ws.Range("F20").Value = 30208
$ws.Range("G20").Value = 781242

$ws.Range("C21").Value = "30656412"
$ws.Range("D21").Value = "BIBIANA LOPEZ DIAZ"
$ws.Range("E21").Value = "1804"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 781242

$wb.Save()
